# Weekly price update for "Fruta, Vega Monumental Concepción - Uva".
# Two new report rows (166 & 167) are inserted into the data table; all
# existing rows from the old row 166 onward shift down by two rows
# (old row 166 -> new row 168, ..., old row 256 -> new row 258).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 166, pushing the existing data (old rows
# 166-256) down to rows 168-258.
$ws.Rows.Item(166).Insert()
$ws.Rows.Item(166).Insert()

# New row 166: Crimpson Seedless / Primera, Región Metropolitana.
$ws.Range("A166").Value = 11
$ws.Range("B166").Value = "Vega Monumental Concepción"
$ws.Range("C166").Value = "Bíobío"
$ws.Range("D166").Value = 45086
$ws.Range("E166").Value = 8
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100109
$ws.Range("H166").Value = "Uva"
$ws.Range("I166").Value = 100109001
$ws.Range("J166").Value = "Uva"
$ws.Range("K166").Value = "Crimpson Seedless"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 100
$ws.Range("N166").Value = 15000
$ws.Range("O166").Value = 16000
$ws.Range("P166").Value = 15500
$ws.Range("Q166").Value = "`$/bandeja 18 kilos"
$ws.Range("R166").Value = "Región Metropolitana"
$ws.Range("S166").Value = 861
$ws.Range("T166").Value = 18

# New row 167: Thompson seedless / Primera, Provincia de San Felipe de Aconcagua.
$ws.Range("A167").Value = 11
$ws.Range("B167").Value = "Vega Monumental Concepción"
$ws.Range("C167").Value = "Bíobío"
$ws.Range("D167").Value = 45086
$ws.Range("E167").Value = 8
$ws.Range("F167").Value = "Fruta"
$ws.Range("G167").Value = 100109
$ws.Range("H167").Value = "Uva"
$ws.Range("I167").Value = 100109001
$ws.Range("J167").Value = "Uva"
$ws.Range("K167").Value = "Thompson seedless"
$ws.Range("L167").Value = "Primera"
$ws.Range("M167").Value = 50
$ws.Range("N167").Value = 30000
$ws.Range("O167").Value = 32000
$ws.Range("P167").Value = 30800
$ws.Range("Q167").Value = "`$/bandeja 18 kilos"
$ws.Range("R167").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S167").Value = 1711
$ws.Range("T167").Value = 18
